# Update Category (column G) values from the old "(AI)" label to the
# rule-derived label that matches the pattern already used by the rest
# of the dataset, per "cashflow and vendor analysis" fix (vendor
# analysis button was producing generic "(AI)" tags instead of the
# correct rule-based category).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Operating Activities (Rule-Default)"
$ws.Range("G5").Value = "Operating Activities (Rule-Vendor)"
$ws.Range("G6").Value = "Operating Activities (Rule-Default)"
$ws.Range("G11").Value = "Operating Activities (Rule-Utility)"
$ws.Range("G20").Value = "Operating Activities (Rule-Default)"
$ws.Range("G21").Value = "Operating Activities (Rule-Default)"
$ws.Range("G24").Value = "Operating Activities (Rule-Payroll)"
$ws.Range("G41").Value = "Financing Activities (Rule-Finance)"
$ws.Range("G47").Value = "Financing Activities (Rule-Finance)"
$ws.Range("G55").Value = "Operating Activities (Rule-Default)"
$ws.Range("G56").Value = "Operating Activities (Rule-Default)"
$ws.Range("G77").Value = "Operating Activities (Rule-Utility)"
$ws.Range("G89").Value = "Operating Activities (Rule-Payroll)"
$ws.Range("G94").Value = "Operating Activities (Rule-Utility)"
$ws.Range("G96").Value = "Operating Activities (Rule-Default)"
$ws.Range("G101").Value = "Operating Activities (Rule-Vendor)"
$ws.Range("G117").Value = "Operating Activities (Rule-Default)"
$ws.Range("G118").Value = "Operating Activities (Rule-Default)"
$ws.Range("G123").Value = "Operating Activities (Rule-Default)"
$ws.Range("G124").Value = "Operating Activities (Rule-Default)"
$ws.Range("G131").Value = "Operating Activities (Rule-Utility)"
$ws.Range("G133").Value = "Operating Activities (Rule-Utility)"
$ws.Range("G134").Value = "Operating Activities (Rule-Vendor)"
$ws.Range("G141").Value = "Operating Activities (Rule-Utility)"
$ws.Range("G147").Value = "Operating Activities (Rule-Utility)"
$ws.Range("G148").Value = "Operating Activities (Rule-Utility)"
$ws.Range("G151").Value = "Operating Activities (Rule-Default)"
$ws.Range("G162").Value = "Operating Activities (Rule-Default)"
$ws.Range("G167").Value = "Financing Activities (Rule-Finance)"
$ws.Range("G168").Value = "Operating Activities (Rule-Default)"
$ws.Range("G169").Value = "Operating Activities (Rule-Default)"
$ws.Range("G183").Value = "Operating Activities (Rule-Vendor)"
$ws.Range("G185").Value = "Operating Activities (Rule-Default)"
$ws.Range("G187").Value = "Operating Activities (Rule-Default)"
$ws.Range("G193").Value = "Operating Activities (Rule-Vendor)"
$ws.Range("G201").Value = "Financing Activities (Rule-Finance)"
$ws.Range("G203").Value = "Operating Activities (Rule-Default)"
$ws.Range("G204").Value = "Operating Activities (Rule-Default)"
$ws.Range("G206").Value = "Operating Activities (Rule-Vendor)"
$ws.Range("G207").Value = "Operating Activities (Rule-Default)"
$ws.Range("G209").Value = "Operating Activities (Rule-Default)"
$ws.Range("G213").Value = "Operating Activities (Rule-Default)"
$ws.Range("G220").Value = "Operating Activities (Rule-Default)"
$ws.Range("G226").Value = "Operating Activities (Rule-Default)"
$ws.Range("G227").Value = "Operating Activities (Rule-Default)"
$ws.Range("G231").Value = "Operating Activities (Rule-Payroll)"
$ws.Range("G233").Value = "Financing Activities (Rule-Finance)"
$ws.Range("G234").Value = "Operating Activities (Rule-Utility)"
$ws.Range("G250").Value = "Operating Activities (Rule-Default)"
$ws.Range("G267").Value = "Operating Activities (Rule-Default)"
$ws.Range("G269").Value = "Operating Activities (Rule-Default)"
$ws.Range("G272").Value = "Operating Activities (Rule-Vendor)"
$ws.Range("G274").Value = "Operating Activities (Rule-Payroll)"
$ws.Range("G275").Value = "Operating Activities (Rule-Utility)"
$ws.Range("G279").Value = "Operating Activities (Rule-Utility)"
$ws.Range("G280").Value = "Financing Activities (Rule-Finance)"
$ws.Range("G283").Value = "Operating Activities (Rule-Default)"
$ws.Range("G290").Value = "Financing Activities (Rule-Finance)"
$ws.Range("G291").Value = "Operating Activities (Rule-Utility)"
$ws.Range("G295").Value = "Operating Activities (Rule-Utility)"
$ws.Range("G301").Value = "Operating Activities (Rule-Utility)"
$ws.Range("G305").Value = "Operating Activities (Rule-Default)"
$ws.Range("G313").Value = "Operating Activities (Rule-Vendor)"
$ws.Range("G314").Value = "Operating Activities (Rule-Utility)"
$ws.Range("G318").Value = "Operating Activities (Rule-Default)"
$ws.Range("G327").Value = "Operating Activities (Rule-Default)"
$ws.Range("G332").Value = "Operating Activities (Rule-Payroll)"
$ws.Range("G334").Value = "Financing Activities (Rule-Finance)"
$ws.Range("G338").Value = "Operating Activities (Rule-Utility)"
$ws.Range("G340").Value = "Operating Activities (Rule-Vendor)"
$ws.Range("G348").Value = "Operating Activities (Rule-Default)"
$ws.Range("G350").Value = "Operating Activities (Rule-Utility)"
$ws.Range("G353").Value = "Operating Activities (Rule-Default)"
$ws.Range("G358").Value = "Operating Activities (Rule-Default)"
$ws.Range("G361").Value = "Operating Activities (Rule-Utility)"
$ws.Range("G373").Value = "Operating Activities (Rule-Default)"
$ws.Range("G377").Value = "Operating Activities (Rule-Default)"
$ws.Range("G379").Value = "Operating Activities (Rule-Vendor)"
$ws.Range("G383").Value = "Operating Activities (Rule-Default)"
$ws.Range("G385").Value = "Operating Activities (Rule-Default)"
$ws.Range("G393").Value = "Financing Activities (Rule-Finance)"
$ws.Range("G399").Value = "Financing Activities (Rule-Finance)"
$ws.Range("G405").Value = "Operating Activities (Rule-Utility)"
$ws.Range("G424").Value = "Financing Activities (Rule-Finance)"
$ws.Range("G425").Value = "Operating Activities (Rule-Default)"
$ws.Range("G427").Value = "Operating Activities (Rule-Vendor)"
$ws.Range("G431").Value = "Financing Activities (Rule-Finance)"
$ws.Range("G436").Value = "Operating Activities (Rule-Default)"
$ws.Range("G439").Value = "Operating Activities (Rule-Utility)"
$ws.Range("G443").Value = "Operating Activities (Rule-Default)"
$ws.Range("G446").Value = "Operating Activities (Rule-Default)"
$ws.Range("G451").Value = "Operating Activities (Rule-Default)"
$ws.Range("G457").Value = "Operating Activities (Rule-Default)"
$ws.Range("G463").Value = "Operating Activities (Rule-Utility)"
$ws.Range("G475").Value = "Operating Activities (Rule-Default)"
$ws.Range("G478").Value = "Operating Activities (Rule-Default)"
$ws.Range("G480").Value = "Operating Activities (Rule-Default)"
$ws.Range("G482").Value = "Operating Activities (Rule-Default)"
$ws.Range("G484").Value = "Operating Activities (Rule-Utility)"
$ws.Range("G491").Value = "Operating Activities (Rule-Default)"

Write-Host "Updated $(($wb.ActiveSheet.UsedRange.Rows.Count)) row sheet: 100 Category cells reclassified from (AI) to Rule-based labels."
